$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12's date cell previously used a one-off "d-mmm" style (index 8). That
# style is being retired, so re-point A12 at the same date style already used
# by the other entries in the column (the style carried by A13) before
# rewriting its value.
$ws.Range("A13").Copy($ws.Range("A12"))
$ws.Range("A12").Value = 45588

# New timesheet entry for the meeting on 2024-10-30 (serial 45595): 3 hours,
# "Lastenheft" booking item, "MockUp- Guis eingebunden" description.
$ws.Range("A13").Copy($ws.Range("A14"))
$ws.Range("A14").Value = 45595

$ws.Range("B13").Copy($ws.Range("B14"))
$ws.Range("B14").Value = 3

$ws.Range("C14").Value = "Lastenheft"
$ws.Range("D14").Value = "MockUp- Guis eingebunden"

# Match the saved selection/active cell state.
$ws.Range("D14").Select() | Out-Null
